$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 5-7 (existing rows whose doctyp_code values changed)
# and add new rows 8-36 for the new document types (Mac-Address / Document Types)

$ws.Range("A5").Value = "DOC001"
$ws.Range("B5").Value = "POI"
$ws.Range("C5").Value = "ara"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "superadmin"
$ws.Range("F5").Value = "now()"

$ws.Range("A6").Value = "CRN"
$ws.Range("B6").Value = "POR"
$ws.Range("C6").Value = "ara"
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "superadmin"
$ws.Range("F6").Value = "now()"

$ws.Range("A7").Value = "COB"
$ws.Range("B7").Value = "POB"
$ws.Range("C7").Value = "ara"
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = "superadmin"
$ws.Range("F7").Value = "now()"

$ws.Range("A8").Value = "DOC001"
$ws.Range("B8").Value = "POI"
$ws.Range("C8").Value = "ara"
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = "superadmin"
$ws.Range("F8").Value = "now()"

$ws.Range("A9").Value = "DOC002"
$ws.Range("B9").Value = "POI"
$ws.Range("C9").Value = "ara"
$ws.Range("D9").Value = $true
$ws.Range("E9").Value = "superadmin"
$ws.Range("F9").Value = "now()"

$ws.Range("A10").Value = "DOC003"
$ws.Range("B10").Value = "POI"
$ws.Range("C10").Value = "ara"
$ws.Range("D10").Value = $true
$ws.Range("E10").Value = "superadmin"
$ws.Range("F10").Value = "now()"

$ws.Range("A11").Value = "DOC004"
$ws.Range("B11").Value = "POI"
$ws.Range("C11").Value = "ara"
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = "superadmin"
$ws.Range("F11").Value = "now()"

$ws.Range("A12").Value = "DOC005"
$ws.Range("B12").Value = "POI"
$ws.Range("C12").Value = "ara"
$ws.Range("D12").Value = $true
$ws.Range("E12").Value = "superadmin"
$ws.Range("F12").Value = "now()"

$ws.Range("A13").Value = "DOC006"
$ws.Range("B13").Value = "POI"
$ws.Range("C13").Value = "ara"
$ws.Range("D13").Value = $true
$ws.Range("E13").Value = "superadmin"
$ws.Range("F13").Value = "now()"

$ws.Range("A14").Value = "DOC007"
$ws.Range("B14").Value = "POI"
$ws.Range("C14").Value = "ara"
$ws.Range("D14").Value = $true
$ws.Range("E14").Value = "superadmin"
$ws.Range("F14").Value = "now()"

$ws.Range("A15").Value = "DOC008"
$ws.Range("B15").Value = "POI"
$ws.Range("C15").Value = "ara"
$ws.Range("D15").Value = $true
$ws.Range("E15").Value = "superadmin"
$ws.Range("F15").Value = "now()"

$ws.Range("A16").Value = "DOC009"
$ws.Range("B16").Value = "POI"
$ws.Range("C16").Value = "ara"
$ws.Range("D16").Value = $true
$ws.Range("E16").Value = "superadmin"
$ws.Range("F16").Value = "now()"

$ws.Range("A17").Value = "DOC010"
$ws.Range("B17").Value = "POI"
$ws.Range("C17").Value = "ara"
$ws.Range("D17").Value = $true
$ws.Range("E17").Value = "superadmin"
$ws.Range("F17").Value = "now()"

$ws.Range("A18").Value = "DOC011"
$ws.Range("B18").Value = "POI"
$ws.Range("C18").Value = "ara"
$ws.Range("D18").Value = $true
$ws.Range("E18").Value = "superadmin"
$ws.Range("F18").Value = "now()"

$ws.Range("A19").Value = "DOC012"
$ws.Range("B19").Value = "POI"
$ws.Range("C19").Value = "ara"
$ws.Range("D19").Value = $true
$ws.Range("E19").Value = "superadmin"
$ws.Range("F19").Value = "now()"

$ws.Range("A20").Value = "DOC001"
$ws.Range("B20").Value = "POA"
$ws.Range("C20").Value = "ara"
$ws.Range("D20").Value = $true
$ws.Range("E20").Value = "superadmin"
$ws.Range("F20").Value = "now()"

$ws.Range("A21").Value = "DOC013"
$ws.Range("B21").Value = "POA"
$ws.Range("C21").Value = "ara"
$ws.Range("D21").Value = $true
$ws.Range("E21").Value = "superadmin"
$ws.Range("F21").Value = "now()"

$ws.Range("A22").Value = "DOC014"
$ws.Range("B22").Value = "POA"
$ws.Range("C22").Value = "ara"
$ws.Range("D22").Value = $true
$ws.Range("E22").Value = "superadmin"
$ws.Range("F22").Value = "now()"

$ws.Range("A23").Value = "DOC015"
$ws.Range("B23").Value = "POA"
$ws.Range("C23").Value = "ara"
$ws.Range("D23").Value = $true
$ws.Range("E23").Value = "superadmin"
$ws.Range("F23").Value = "now()"

$ws.Range("A24").Value = "DOC004"
$ws.Range("B24").Value = "POA"
$ws.Range("C24").Value = "ara"
$ws.Range("D24").Value = $true
$ws.Range("E24").Value = "superadmin"
$ws.Range("F24").Value = "now()"

$ws.Range("A25").Value = "DOC005"
$ws.Range("B25").Value = "POA"
$ws.Range("C25").Value = "ara"
$ws.Range("D25").Value = $true
$ws.Range("E25").Value = "superadmin"
$ws.Range("F25").Value = "now()"

$ws.Range("A26").Value = "DOC006"
$ws.Range("B26").Value = "POA"
$ws.Range("C26").Value = "ara"
$ws.Range("D26").Value = $true
$ws.Range("E26").Value = "superadmin"
$ws.Range("F26").Value = "now()"

$ws.Range("A27").Value = "DOC016"
$ws.Range("B27").Value = "POA"
$ws.Range("C27").Value = "ara"
$ws.Range("D27").Value = $true
$ws.Range("E27").Value = "superadmin"
$ws.Range("F27").Value = "now()"

$ws.Range("A28").Value = "DOC017"
$ws.Range("B28").Value = "POA"
$ws.Range("C28").Value = "ara"
$ws.Range("D28").Value = $true
$ws.Range("E28").Value = "superadmin"
$ws.Range("F28").Value = "now()"

$ws.Range("A29").Value = "DOC018"
$ws.Range("B29").Value = "POA"
$ws.Range("C29").Value = "ara"
$ws.Range("D29").Value = $true
$ws.Range("E29").Value = "superadmin"
$ws.Range("F29").Value = "now()"

$ws.Range("A30").Value = "DOC008"
$ws.Range("B30").Value = "POA"
$ws.Range("C30").Value = "ara"
$ws.Range("D30").Value = $true
$ws.Range("E30").Value = "superadmin"
$ws.Range("F30").Value = "now()"

$ws.Range("A31").Value = "DOC024"
$ws.Range("B31").Value = "POR"
$ws.Range("C31").Value = "ara"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = "DOC025"
$ws.Range("B32").Value = "POR"
$ws.Range("C32").Value = "ara"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

$ws.Range("A33").Value = "DOC026"
$ws.Range("B33").Value = "POR"
$ws.Range("C33").Value = "ara"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

$ws.Range("A34").Value = "DOC001"
$ws.Range("B34").Value = "POR"
$ws.Range("C34").Value = "ara"
$ws.Range("D34").Value = $true
$ws.Range("E34").Value = "superadmin"
$ws.Range("F34").Value = "now()"

$ws.Range("A35").Value = "DOC027"
$ws.Range("B35").Value = "POR"
$ws.Range("C35").Value = "ara"
$ws.Range("D35").Value = $true
$ws.Range("E35").Value = "superadmin"
$ws.Range("F35").Value = "now()"

$ws.Range("A36").Value = "DOC028"
$ws.Range("B36").Value = "POR"
$ws.Range("C36").Value = "ara"
$ws.Range("D36").Value = $true
$ws.Range("E36").Value = "superadmin"
$ws.Range("F36").Value = "now()"

# Update the visible selection to mirror the authored state (G1:XFD1048576)
$ws.Range("G1:XFD1048576").Select()

